$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("E2").Value = 15; $ws.Range("F2").Value = 14; $ws.Range("G2").Value = 51.72; $ws.Range("H2").Value = 48.28; $ws.Range("I2").Value = 6.7; $ws.Range("J2").Value = 0; $ws.Range("K2").Value = 0
$ws.Range("E3").Value = 10; $ws.Range("F3").Value = 11; $ws.Range("G3").Value = 47.62; $ws.Range("H3").Value = 52.38; $ws.Range("I3").Value = 6.2; $ws.Range("J3").Value = 0; $ws.Range("K3").Value = 0
$ws.Range("E4").Value = 21; $ws.Range("F4").Value = 7; $ws.Range("G4").Value = 75; $ws.Range("H4").Value = 25; $ws.Range("I4").Value = 7.3; $ws.Range("J4").Value = 0; $ws.Range("K4").Value = 0
$ws.Range("I5").Value = 5.7; $ws.Range("J5").Value = 0; $ws.Range("K5").Value = 0
$ws.Range("E6").Value = 7; $ws.Range("F6").Value = 15; $ws.Range("G6").Value = 31.82; $ws.Range("H6").Value = 68.18000000000001; $ws.Range("I6").Value = 5.5; $ws.Range("J6").Value = 0; $ws.Range("K6").Value = 0
$ws.Range("E7").Value = 34; $ws.Range("F7").Value = 5; $ws.Range("G7").Value = 87.18000000000001; $ws.Range("H7").Value = 12.82; $ws.Range("I7").Value = 6.7; $ws.Range("J7").Value = 3; $ws.Range("K7").Value = 7.69
$ws.Range("E8").Value = 38; $ws.Range("F8").Value = 3; $ws.Range("G8").Value = 92.68000000000001; $ws.Range("H8").Value = 7.32; $ws.Range("I8").Value = 7.2; $ws.Range("J8").Value = 0; $ws.Range("K8").Value = 0
$ws.Range("E9").Value = 24; $ws.Range("F9").Value = 1; $ws.Range("G9").Value = 96; $ws.Range("H9").Value = 4; $ws.Range("I9").Value = 6.6; $ws.Range("J9").Value = 0; $ws.Range("K9").Value = 0
$ws.Range("E10").Value = 36; $ws.Range("F10").Value = 3; $ws.Range("G10").Value = 92.31; $ws.Range("H10").Value = 7.69; $ws.Range("I10").Value = 7.2; $ws.Range("J10").Value = 2; $ws.Range("K10").Value = 5.13
$ws.Range("I11").Value = 7.2; $ws.Range("J11").Value = 2; $ws.Range("K11").Value = 5.71
$ws.Range("E12").Value = 29; $ws.Range("F12").Value = 7; $ws.Range("G12").Value = 80.56; $ws.Range("H12").Value = 19.44; $ws.Range("I12").Value = 6.6; $ws.Range("J12").Value = 1; $ws.Range("K12").Value = 2.78
$ws.Range("E13").Value = 18; $ws.Range("F13").Value = 18; $ws.Range("G13").Value = 50; $ws.Range("H13").Value = 50; $ws.Range("I13").Value = 6.2; $ws.Range("J13").Value = 11; $ws.Range("K13").Value = 30.56
$ws.Range("E18").Value = 30; $ws.Range("F18").Value = 3; $ws.Range("G18").Value = 90.91; $ws.Range("H18").Value = 9.09; $ws.Range("I18").Value = 8; $ws.Range("J18").Value = 3; $ws.Range("K18").Value = 9.09
$ws.Range("E19").Value = 21; $ws.Range("F19").Value = 17; $ws.Range("G19").Value = 55.26; $ws.Range("H19").Value = 44.74; $ws.Range("I19").Value = 8.4; $ws.Range("J19").Value = 17; $ws.Range("K19").Value = 44.74
$ws.Range("E20").Value = 29; $ws.Range("F20").Value = 10; $ws.Range("G20").Value = 74.36; $ws.Range("H20").Value = 25.64; $ws.Range("I20").Value = 8.199999999999999; $ws.Range("J20").Value = 10; $ws.Range("K20").Value = 25.64
$ws.Range("E21").Value = 22; $ws.Range("F21").Value = 12; $ws.Range("G21").Value = 64.70999999999999; $ws.Range("H21").Value = 35.29; $ws.Range("I21").Value = 8.800000000000001; $ws.Range("J21").Value = 12; $ws.Range("K21").Value = 35.29
$ws.Range("E28").Value = 18; $ws.Range("F28").Value = 21; $ws.Range("G28").Value = 46.15; $ws.Range("H28").Value = 53.85; $ws.Range("I28").Value = 6.4; $ws.Range("J28").Value = 12; $ws.Range("K28").Value = 30.77
$ws.Range("E30").Value = 15; $ws.Range("F30").Value = 9; $ws.Range("G30").Value = 62.5; $ws.Range("H30").Value = 37.5; $ws.Range("I30").Value = 6.4; $ws.Range("J30").Value = 1; $ws.Range("K30").Value = 4.17
$ws.Range("E31").Value = 15; $ws.Range("F31").Value = 22; $ws.Range("G31").Value = 40.54; $ws.Range("H31").Value = 59.46; $ws.Range("I31").Value = 5.8; $ws.Range("J31").Value = 0; $ws.Range("K31").Value = 0

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("E2").Value = 16; $ws.Range("F2").Value = 13; $ws.Range("G2").Value = 55.17; $ws.Range("H2").Value = 44.83; $ws.Range("I2").Value = 6.7; $ws.Range("J2").Value = 0; $ws.Range("K2").Value = 0
$ws.Range("E3").Value = 10; $ws.Range("F3").Value = 11; $ws.Range("G3").Value = 47.62; $ws.Range("H3").Value = 52.38; $ws.Range("I3").Value = 5.7; $ws.Range("J3").Value = 0; $ws.Range("K3").Value = 0
$ws.Range("E4").Value = 23; $ws.Range("F4").Value = 5; $ws.Range("G4").Value = 82.14; $ws.Range("H4").Value = 17.86; $ws.Range("I4").Value = 7.1; $ws.Range("J4").Value = 0; $ws.Range("K4").Value = 0
$ws.Range("E5").Value = 8; $ws.Range("F5").Value = 14; $ws.Range("G5").Value = 36.36; $ws.Range("H5").Value = 63.64; $ws.Range("I5").Value = 6.4; $ws.Range("J5").Value = 0; $ws.Range("K5").Value = 0
$ws.Range("E6").Value = 7; $ws.Range("F6").Value = 15; $ws.Range("G6").Value = 31.82; $ws.Range("H6").Value = 68.18000000000001; $ws.Range("I6").Value = 6.3; $ws.Range("J6").Value = 0; $ws.Range("K6").Value = 0
$ws.Range("E7").Value = 25; $ws.Range("F7").Value = 14; $ws.Range("G7").Value = 64.09999999999999; $ws.Range("H7").Value = 35.9; $ws.Range("I7").Value = 6.7; $ws.Range("J7").Value = 14; $ws.Range("K7").Value = 35.9
$ws.Range("E8").Value = 34; $ws.Range("F8").Value = 7; $ws.Range("G8").Value = 82.93000000000001; $ws.Range("H8").Value = 17.07; $ws.Range("I8").Value = 7; $ws.Range("J8").Value = 7; $ws.Range("K8").Value = 17.07
$ws.Range("E9").Value = 17; $ws.Range("F9").Value = 8; $ws.Range("G9").Value = 68; $ws.Range("H9").Value = 32; $ws.Range("I9").Value = 6.9; $ws.Range("J9").Value = 8; $ws.Range("K9").Value = 32
$ws.Range("E10").Value = 28; $ws.Range("F10").Value = 11; $ws.Range("G10").Value = 71.79000000000001; $ws.Range("H10").Value = 28.21; $ws.Range("I10").Value = 7.1; $ws.Range("J10").Value = 11; $ws.Range("K10").Value = 28.21
$ws.Range("E11").Value = 22; $ws.Range("F11").Value = 13; $ws.Range("G11").Value = 62.86; $ws.Range("H11").Value = 37.14; $ws.Range("I11").Value = 6.8; $ws.Range("J11").Value = 13; $ws.Range("K11").Value = 37.14
$ws.Range("E12").Value = 21; $ws.Range("F12").Value = 15; $ws.Range("G12").Value = 58.33; $ws.Range("H12").Value = 41.67; $ws.Range("I12").Value = 7; $ws.Range("J12").Value = 14; $ws.Range("K12").Value = 38.89
$ws.Range("E13").Value = 13; $ws.Range("F13").Value = 23; $ws.Range("G13").Value = 36.11; $ws.Range("H13").Value = 63.89; $ws.Range("I13").Value = 6.7; $ws.Range("J13").Value = 20; $ws.Range("K13").Value = 55.56
$ws.Range("E18").Value = 24; $ws.Range("F18").Value = 9; $ws.Range("G18").Value = 72.73; $ws.Range("H18").Value = 27.27; $ws.Range("I18").Value = 7.9; $ws.Range("J18").Value = 9; $ws.Range("K18").Value = 27.27
$ws.Range("E19").Value = 17; $ws.Range("F19").Value = 21; $ws.Range("G19").Value = 44.74; $ws.Range("H19").Value = 55.26; $ws.Range("I19").Value = 8; $ws.Range("J19").Value = 21; $ws.Range("K19").Value = 55.26
$ws.Range("E20").Value = 22; $ws.Range("F20").Value = 17; $ws.Range("G20").Value = 56.41; $ws.Range("H20").Value = 43.59; $ws.Range("I20").Value = 8; $ws.Range("J20").Value = 17; $ws.Range("K20").Value = 43.59
$ws.Range("E21").Value = 20; $ws.Range("F21").Value = 14; $ws.Range("G21").Value = 58.82; $ws.Range("H21").Value = 41.18; $ws.Range("I21").Value = 8.199999999999999; $ws.Range("J21").Value = 14; $ws.Range("K21").Value = 41.18
$ws.Range("E28").Value = 13; $ws.Range("F28").Value = 26; $ws.Range("G28").Value = 33.33; $ws.Range("H28").Value = 66.67; $ws.Range("I28").Value = 6.5; $ws.Range("J28").Value = 19; $ws.Range("K28").Value = 48.72
$ws.Range("E30").Value = 11; $ws.Range("F30").Value = 13; $ws.Range("G30").Value = 45.83; $ws.Range("H30").Value = 54.17; $ws.Range("I30").Value = 6.3; $ws.Range("J30").Value = 6; $ws.Range("K30").Value = 25
$ws.Range("E31").Value = 4; $ws.Range("F31").Value = 33; $ws.Range("G31").Value = 10.81; $ws.Range("H31").Value = 89.19; $ws.Range("I31").Value = 6.8; $ws.Range("J31").Value = 32; $ws.Range("K31").Value = 86.48999999999999

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("E2").Value = 16; $ws.Range("F2").Value = 13; $ws.Range("G2").Value = 55.17; $ws.Range("H2").Value = 44.83; $ws.Range("I2").Value = 6.8; $ws.Range("J2").Value = 0; $ws.Range("K2").Value = 0
$ws.Range("E3").Value = 10; $ws.Range("F3").Value = 11; $ws.Range("G3").Value = 47.62; $ws.Range("H3").Value = 52.38; $ws.Range("I3").Value = 6; $ws.Range("J3").Value = 0; $ws.Range("K3").Value = 0
$ws.Range("E4").Value = 23; $ws.Range("F4").Value = 5; $ws.Range("G4").Value = 82.14; $ws.Range("H4").Value = 17.86; $ws.Range("I4").Value = 7.4; $ws.Range("J4").Value = 0; $ws.Range("K4").Value = 0
$ws.Range("E5").Value = 8; $ws.Range("F5").Value = 14; $ws.Range("G5").Value = 36.36; $ws.Range("H5").Value = 63.64; $ws.Range("I5").Value = 6; $ws.Range("J5").Value = 0; $ws.Range("K5").Value = 0
$ws.Range("E6").Value = 7; $ws.Range("F6").Value = 15; $ws.Range("G6").Value = 31.82; $ws.Range("H6").Value = 68.18000000000001; $ws.Range("I6").Value = 5.8; $ws.Range("J6").Value = 0; $ws.Range("K6").Value = 0
$ws.Range("E7").Value = 34; $ws.Range("F7").Value = 5; $ws.Range("G7").Value = 87.18000000000001; $ws.Range("H7").Value = 12.82; $ws.Range("I7").Value = 6.7; $ws.Range("J7").Value = 3; $ws.Range("K7").Value = 7.69
$ws.Range("E8").Value = 38; $ws.Range("F8").Value = 3; $ws.Range("G8").Value = 92.68000000000001; $ws.Range("H8").Value = 7.32; $ws.Range("I8").Value = 7.2; $ws.Range("J8").Value = 0; $ws.Range("K8").Value = 0
$ws.Range("E9").Value = 24; $ws.Range("F9").Value = 1; $ws.Range("G9").Value = 96; $ws.Range("H9").Value = 4; $ws.Range("I9").Value = 6.7; $ws.Range("J9").Value = 0; $ws.Range("K9").Value = 0
$ws.Range("E10").Value = 36; $ws.Range("F10").Value = 3; $ws.Range("G10").Value = 92.31; $ws.Range("H10").Value = 7.69; $ws.Range("I10").Value = 7.2; $ws.Range("J10").Value = 2; $ws.Range("K10").Value = 5.13
$ws.Range("I11").Value = 7.2; $ws.Range("J11").Value = 2; $ws.Range("K11").Value = 5.71
$ws.Range("E12").Value = 29; $ws.Range("F12").Value = 7; $ws.Range("G12").Value = 80.56; $ws.Range("H12").Value = 19.44; $ws.Range("I12").Value = 6.7; $ws.Range("J12").Value = 1; $ws.Range("K12").Value = 2.78
$ws.Range("E13").Value = 18; $ws.Range("F13").Value = 18; $ws.Range("G13").Value = 50; $ws.Range("H13").Value = 50; $ws.Range("I13").Value = 6.4; $ws.Range("J13").Value = 11; $ws.Range("K13").Value = 30.56
$ws.Range("E18").Value = 30; $ws.Range("F18").Value = 3; $ws.Range("G18").Value = 90.91; $ws.Range("H18").Value = 9.09; $ws.Range("J18").Value = 3; $ws.Range("K18").Value = 9.09
$ws.Range("E19").Value = 21; $ws.Range("F19").Value = 17; $ws.Range("G19").Value = 55.26; $ws.Range("H19").Value = 44.74; $ws.Range("I19").Value = 8.199999999999999; $ws.Range("J19").Value = 17; $ws.Range("K19").Value = 44.74
$ws.Range("E20").Value = 29; $ws.Range("F20").Value = 10; $ws.Range("G20").Value = 74.36; $ws.Range("H20").Value = 25.64; $ws.Range("I20").Value = 8.199999999999999; $ws.Range("J20").Value = 10; $ws.Range("K20").Value = 25.64
$ws.Range("E21").Value = 22; $ws.Range("F21").Value = 12; $ws.Range("G21").Value = 64.70999999999999; $ws.Range("H21").Value = 35.29; $ws.Range("I21").Value = 8.699999999999999; $ws.Range("J21").Value = 12; $ws.Range("K21").Value = 35.29
$ws.Range("E28").Value = 14; $ws.Range("F28").Value = 25; $ws.Range("G28").Value = 35.9; $ws.Range("H28").Value = 64.09999999999999; $ws.Range("I28").Value = 6.3; $ws.Range("J28").Value = 12; $ws.Range("K28").Value = 30.77
$ws.Range("E30").Value = 14; $ws.Range("F30").Value = 10; $ws.Range("G30").Value = 58.33; $ws.Range("H30").Value = 41.67; $ws.Range("I30").Value = 6.3; $ws.Range("J30").Value = 1; $ws.Range("K30").Value = 4.17
$ws.Range("E31").Value = 14; $ws.Range("F31").Value = 23; $ws.Range("G31").Value = 37.84; $ws.Range("H31").Value = 62.16; $ws.Range("I31").Value = 5.8; $ws.Range("J31").Value = 0; $ws.Range("K31").Value = 0
